# Adapt the column header formatting to the respective input file names:
#   "<field>_old" -> "<field>_FV2310"
#   "<field>_new" -> "<field>_FV2404"
# Also (re)create the Excel Table over the data range and freeze the header row,
# matching the regenerated xlsx export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the "_old" / "_new" header suffixes to the format-version suffixes ---

$fv2310Headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

# Columns A..J hold the "_old"-suffixed headers (-> "_FV2310")
for ($i = 0; $i -lt $fv2310Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $fv2310Headers[$i]
}

# Column K ("diff") is unchanged.

# Columns L..U hold the "_new"-suffixed headers (-> "_FV2404")
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value2 = $fv2404Headers[$i]
}

# --- 2. (Re)create the worksheet Table over the full data range ---

$dataRange = $ws.Range("A1:U92")
$table = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$table.Name = "Table1"

# --- 3. Freeze the header row ---

$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "edit applied"
